$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column D keeps its text representation (avoid Excel
# auto-converting numeric-looking strings like "1.001" into floating
# point numbers, which would lose/round the original formatting).
$ws.Range("D2:D51").NumberFormat = "@"

$data = @(
    ,@(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "30.823.24", "  -0.29%  ")
    ,@(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.939.11", "  -0.38%  ")
    ,@(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.001", "  +0.41%  ")
    ,@(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "243.89", "  -0.33%  ")
    ,@(6, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.001", "  +0.50%  ")
    ,@(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.4881", "  +0.27%  ")
    ,@(8, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.2950", "  -0.28%  ")
    ,@(9, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.06898", "  +1.17%  ")
    ,@(10, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "19.33", "  +1.06%  ")
    ,@(11, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "105.00", "  -1.78%  ")
    ,@(12, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.07795", "  +0.99%  ")
    ,@(13, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.935.00", "  -0.66%  ")
    ,@(14, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.356", "  -1.90%  ")
    ,@(15, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.7024", "  -0.34%  ")
    ,@(16, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "273.39", "  -2.78%  ")
    ,@(17, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "30.827.15", "  -0.38%  ")
    ,@(18, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000007738", "  +0.16%  ")
    ,@(19, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "13.15", "  -0.56%  ")
    ,@(20, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "5.639", "  +2.68%  ")
    ,@(21, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.001", "  +0.50%  ")
    ,@(22, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.001", "  +0.49%  ")
    ,@(23, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "6.540", "  +0.88%  ")
    ,@(24, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "9.829", "  +0.22%  ")
    ,@(25, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "165.26", "  -2.18%  ")
    ,@(26, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "19.64", "  -1.48%  ")
    ,@(27, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "2.165", "  -2.08%  ")
    ,@(28, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.1039", "  -1.34%  ")
    ,@(29, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.387", "  -1.61%  ")
    ,@(30, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.649", "  +2.12%  ")
    ,@(31, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.559", "  -1.36%  ")
    ,@(32, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "4.422", "  -1.16%  ")
    ,@(33, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.04906", "  -0.84%  ")
    ,@(34, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.7614", "  -0.27%  ")
    ,@(35, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.151", "  -2.07%  ")
    ,@(36, "Frax", "https://coinranking.com/coin/KfWtaeV1W+frax-frax", "1.000", "  +0.46%  ")
    ,@(37, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.734", "  +0.62%  ")
    ,@(38, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.02014", "  -0.48%  ")
    ,@(39, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "79.72", "  +6.18%  ")
    ,@(40, "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.666", "  -0.85%  ")
    ,@(41, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "6.528", "  +0.45%  ")
    ,@(42, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "2.088", "  -3.07%  ")
    ,@(43, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.9052", "  +2.67%  ")
    ,@(44, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.4457", "  -0.65%  ")
    ,@(45, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "108.54", "  -0.60%  ")
    ,@(46, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "7.844", "  -3.74%  ")
    ,@(47, "PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "1.001", "  +0.52%  ")
    ,@(48, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "996.94", "  +1.81%  ")
    ,@(49, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.1251", "  -0.84%  ")
    ,@(50, "Elrond", "https://coinranking.com/coin/omwkOTglq+elrond-egld", "36.31", "  +1.68%  ")
    ,@(51, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "9.220", "  -1.65%  ")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
